$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.771.79'
$ws.Range('E2').Value = '  +0.11%  '
$ws.Range('D3').Value = '1.628.84'
$ws.Range('E3').Value = '  -0.35%  '
$ws.Range('E4').Value = '  -0.64%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.98'
$ws.Range('E5').Value = '  -0.65%  '
$ws.Range('E6').Value = '  -0.30%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.49%  '
$ws.Range('E8').Value = '  -1.01%  '
$ws.Range('E9').Value = '  -0.98%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.58'
$ws.Range('E10').Value = '  +0.02%  '
$ws.Range('E11').Value = '  +0.24%  '
$ws.Range('E12').Value = '  +0.31%  '
$ws.Range('D13').Value = '1.853.86'
$ws.Range('E13').Value = '  -0.32%  '
$ws.Range('D14').Value = '1.617.49'
$ws.Range('E14').Value = '  -1.11%  '
$ws.Range('E15').Value = '  -1.03%  '
$ws.Range('D16').Value = '0.0₃0756'
$ws.Range('E16').Value = '  -1.28%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.60'
$ws.Range('E17').Value = '  -0.14%  '
$ws.Range('D18').Value = '25.791.11'
$ws.Range('E18').Value = '  +0.16%  '
$ws.Range('E19').Value = '  -0.54%  '
$ws.Range('E20').Value = '  -0.23%  '
$ws.Range('E21').Value = '  -1.56%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.92'
$ws.Range('E22').Value = '  -0.22%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.28'
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.998'
$ws.Range('E24').Value = '  -0.52%  '
$ws.Range('E25').Value = '  -2.32%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '142.26'
$ws.Range('E26').Value = '  +1.51%  '
$ws.Range('E27').Value = '  +0.57%  '
$ws.Range('E28').Value = '  -0.83%  '
$ws.Range('E29').Value = '  -0.12%  '
$ws.Range('E30').Value = '  -0.68%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0494'
$ws.Range('E31').Value = '  +0.09%  '
$ws.Range('E32').Value = '  -0.62%  '
$ws.Range('E33').Value = '  -0.81%  '
$ws.Range('E34').Value = '  -0.37%  '
$ws.Range('E35').Value = '  +0.12%  '
$ws.Range('E36').Value = '  +0.37%  '
$ws.Range('D37').Value = '1.141.70'
$ws.Range('E37').Value = '  +2.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.543'
$ws.Range('E39').Value = '  -1.35%  '
$ws.Range('E40').Value = '  -0.25%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.997'
$ws.Range('E41').Value = '  -0.61%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.61'
$ws.Range('E42').Value = '  +0.99%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.43'
$ws.Range('E43').Value = '  +0.74%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.800'
$ws.Range('E44').Value = '  -0.27%  '
$ws.Range('D45').Value = '1.765.03'
$ws.Range('E45').Value = '  -0.22%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '55.35'
$ws.Range('E46').Value = '  +0.39%  '
$ws.Range('E47').Value = '  +2.15%  '
$ws.Range('E48').Value = '  +5.32%  '
$ws.Range('E49').Value = '  -0.39%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.55'
$ws.Range('E50').Value = '  -0.48%  '
$ws.Range('E51').Value = '  +1.57%  '
